# Reorder the "Recorded By" (column G) comma-separated list so that the
# literal entry "System" (exact case) moves from the end of the list to the
# front, and the rest of the entries are reversed to match the source
# repository's canonical ordering.
#
# Rule (derived from the target diff): for any cell in column G whose value
# is a comma-separated list that contains "System" as one of its elements
# (and has more than one element), the whole list is reversed.  Cells that
# do not contain "System" as a distinct element, or that contain only a
# single element, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $text = $cell.Text

    if ($null -eq $text) { continue }
    if ($text -eq "") { continue }

    $parts = $text -split ",\s*"

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p -ceq "System") { $hasSystem = $true }
    }

    if ($hasSystem -and $parts.Count -gt 1) {
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversed)
        $cell.Value = $newText
    }
}

